# "Reservado arquivo para cfarl"
# For rows 527-557 (excluding 536, which already has a D value), mark column D
# with "cfarl" to flag the file as reserved for that translator.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToFlag = @(527,528,529,530,531,532,533,534,535,537,538,539,540,541,542,543,544,545,546,547,548,549,550,551,552,553,554,555,556,557)
foreach ($r in $rowsToFlag) {
    $ws.Cells.Item($r, 4).Value = "cfarl"
}

# Rows 596-606: mark as translated ("SIM") in column C.
for ($r = 596; $r -le 606; $r++) {
    $ws.Cells.Item($r, 3).Value = "SIM"
}

# Move the current selection/scroll position to reflect where the user was working.
[void]$ws.Range("E530").Select()
$excel.ActiveWindow.ScrollRow = 525
